$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data rows 117 and 118 (id 115 and 116) had their match details swapped:
# column B (match id number) and columns F through AC (HomeTeam, AwayTeam,
# result and all the odds columns) move to the other row, while columns
# A (id), C (Div), D (Div Original Name) and E (Date) stay put.
#
# NOTE: use .Value2 (not .Value) - this runtime's .Value getter does not
# return usable data for Range objects.

$b117 = $ws.Range("B117").Value2
$b118 = $ws.Range("B118").Value2

$block117 = $ws.Range("F117:AC117").Value2
$block118 = $ws.Range("F118:AC118").Value2

$ws.Range("B117").Value2 = $b118
$ws.Range("B118").Value2 = $b117

$ws.Range("F117:AC117").Value2 = $block118
$ws.Range("F118:AC118").Value2 = $block117
